$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at AK (column 37) so the new header cell picks up
# the same formatting as the rest of the header row (mirrors what Excel
# does when a user inserts a column next to existing data).
$ws.Columns.Item(37).Insert()

# The insert also stamps row 2 of the new column with the neighbouring
# format; clear it back out since the source data only adds a header.
$ws.Cells.Item(2, 37).Clear()

# Add the new header "otherAdvisor" in column AK (37), row 1
$ws.Cells.Item(1, 37).Value = "otherAdvisor"

# Match the column width Excel computed for the new column (ColumnWidth
# units round-trip slightly differently than the raw stored width, so
# feed in the value that lands closest to the target ~16.664)
$ws.Columns.Item(37).ColumnWidth = 15.83
